$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45993
$ws.Cells.Item(2, 2).Value = 13
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = "02.12.202513"
$ws.Cells.Item(3, 1).Value = 45993
$ws.Cells.Item(3, 2).Value = 14
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = "02.12.202514"
$ws.Cells.Item(4, 1).Value = 45993
$ws.Cells.Item(4, 2).Value = 15
$ws.Cells.Item(4, 3).Value = 0.365
$ws.Cells.Item(4, 4).Value = "02.12.202515"
$ws.Cells.Item(5, 1).Value = 45993
$ws.Cells.Item(5, 2).Value = 16
$ws.Cells.Item(5, 3).Value = 0.312
$ws.Cells.Item(5, 4).Value = "02.12.202516"
$ws.Cells.Item(6, 1).Value = 45993
$ws.Cells.Item(6, 2).Value = 17
$ws.Cells.Item(6, 3).Value = 0.119
$ws.Cells.Item(6, 4).Value = "02.12.202517"
$ws.Cells.Item(7, 1).Value = 45993
$ws.Cells.Item(7, 2).Value = 18
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = "02.12.202518"
$ws.Cells.Item(8, 1).Value = 45993
$ws.Cells.Item(8, 2).Value = 19
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = "02.12.202519"
$ws.Cells.Item(9, 1).Value = 45993
$ws.Cells.Item(9, 2).Value = 20
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = "02.12.202520"
$ws.Cells.Item(10, 1).Value = 45993
$ws.Cells.Item(10, 2).Value = 21
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = "02.12.202521"
$ws.Cells.Item(11, 1).Value = 45993
$ws.Cells.Item(11, 2).Value = 22
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = "02.12.202522"
$ws.Cells.Item(12, 1).Value = 45993
$ws.Cells.Item(12, 2).Value = 23
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = "02.12.202523"
$ws.Cells.Item(13, 1).Value = 45993
$ws.Cells.Item(13, 2).Value = 24
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = "02.12.202524"
$ws.Cells.Item(14, 1).Value = 45994
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = "03.12.20251"
$ws.Cells.Item(15, 1).Value = 45994
$ws.Cells.Item(15, 2).Value = 2
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = "03.12.20252"
$ws.Cells.Item(16, 1).Value = 45994
$ws.Cells.Item(16, 2).Value = 3
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = "03.12.20253"
$ws.Cells.Item(17, 1).Value = 45994
$ws.Cells.Item(17, 2).Value = 4
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = "03.12.20254"
$ws.Cells.Item(18, 1).Value = 45994
$ws.Cells.Item(18, 2).Value = 5
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = "03.12.20255"
$ws.Cells.Item(19, 1).Value = 45994
$ws.Cells.Item(19, 2).Value = 6
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = "03.12.20256"
$ws.Cells.Item(20, 1).Value = 45994
$ws.Cells.Item(20, 2).Value = 7
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = "03.12.20257"
$ws.Cells.Item(21, 1).Value = 45994
$ws.Cells.Item(21, 2).Value = 8
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 4).Value = "03.12.20258"
$ws.Cells.Item(22, 1).Value = 45994
$ws.Cells.Item(22, 2).Value = 9
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = "03.12.20259"
$ws.Cells.Item(23, 1).Value = 45994
$ws.Cells.Item(23, 2).Value = 10
$ws.Cells.Item(23, 3).Value = 0.083
$ws.Cells.Item(23, 4).Value = "03.12.202510"
$ws.Cells.Item(24, 1).Value = 45994
$ws.Cells.Item(24, 2).Value = 11
$ws.Cells.Item(24, 3).Value = 0.406
$ws.Cells.Item(24, 4).Value = "03.12.202511"
$ws.Cells.Item(25, 1).Value = 45994
$ws.Cells.Item(25, 2).Value = 12
$ws.Cells.Item(25, 3).Value = 0.87
$ws.Cells.Item(25, 4).Value = "03.12.202512"
$ws.Cells.Item(26, 1).Value = 45994
$ws.Cells.Item(26, 2).Value = 13
$ws.Cells.Item(26, 3).Value = 1.223
$ws.Cells.Item(26, 4).Value = "03.12.202513"
$ws.Cells.Item(27, 1).Value = 45994
$ws.Cells.Item(27, 2).Value = 14
$ws.Cells.Item(27, 3).Value = 1.277
$ws.Cells.Item(27, 4).Value = "03.12.202514"
$ws.Cells.Item(28, 1).Value = 45994
$ws.Cells.Item(28, 2).Value = 15
$ws.Cells.Item(28, 3).Value = 1.182
$ws.Cells.Item(28, 4).Value = "03.12.202515"
$ws.Cells.Item(29, 1).Value = 45994
$ws.Cells.Item(29, 2).Value = 16
$ws.Cells.Item(29, 3).Value = 0.767
$ws.Cells.Item(29, 4).Value = "03.12.202516"
$ws.Cells.Item(30, 1).Value = 45994
$ws.Cells.Item(30, 2).Value = 17
$ws.Cells.Item(30, 3).Value = 0.178
$ws.Cells.Item(30, 4).Value = "03.12.202517"
$ws.Cells.Item(31, 1).Value = 45994
$ws.Cells.Item(31, 2).Value = 18
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = "03.12.202518"
$ws.Cells.Item(32, 1).Value = 45994
$ws.Cells.Item(32, 2).Value = 19
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = "03.12.202519"
$ws.Cells.Item(33, 1).Value = 45994
$ws.Cells.Item(33, 2).Value = 20
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = "03.12.202520"
$ws.Cells.Item(34, 1).Value = 45994
$ws.Cells.Item(34, 2).Value = 21
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = "03.12.202521"
$ws.Cells.Item(35, 1).Value = 45994
$ws.Cells.Item(35, 2).Value = 22
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = "03.12.202522"
$ws.Cells.Item(36, 1).Value = 45994
$ws.Cells.Item(36, 2).Value = 23
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 4).Value = "03.12.202523"
$ws.Cells.Item(37, 1).Value = 45994
$ws.Cells.Item(37, 2).Value = 24
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = "03.12.202524"
$ws.Cells.Item(38, 1).Value = 45995
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(38, 4).Value = "04.12.20251"
$ws.Cells.Item(39, 1).Value = 45995
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 4).Value = "04.12.20252"
$ws.Cells.Item(40, 1).Value = 45995
$ws.Cells.Item(40, 2).Value = 3
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = "04.12.20253"
$ws.Cells.Item(41, 1).Value = 45995
$ws.Cells.Item(41, 2).Value = 4
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(41, 4).Value = "04.12.20254"
$ws.Cells.Item(42, 1).Value = 45995
$ws.Cells.Item(42, 2).Value = 5
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(42, 4).Value = "04.12.20255"
$ws.Cells.Item(43, 1).Value = 45995
$ws.Cells.Item(43, 2).Value = 6
$ws.Cells.Item(43, 3).Value = 0
$ws.Cells.Item(43, 4).Value = "04.12.20256"
$ws.Cells.Item(44, 1).Value = 45995
$ws.Cells.Item(44, 2).Value = 7
$ws.Cells.Item(44, 3).Value = 0
$ws.Cells.Item(44, 4).Value = "04.12.20257"
$ws.Cells.Item(45, 1).Value = 45995
$ws.Cells.Item(45, 2).Value = 8
$ws.Cells.Item(45, 3).Value = 0
$ws.Cells.Item(45, 4).Value = "04.12.20258"
$ws.Cells.Item(46, 1).Value = 45995
$ws.Cells.Item(46, 2).Value = 9
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = "04.12.20259"
$ws.Cells.Item(47, 1).Value = 45995
$ws.Cells.Item(47, 2).Value = 10
$ws.Cells.Item(47, 3).Value = 0.083
$ws.Cells.Item(47, 4).Value = "04.12.202510"
$ws.Cells.Item(48, 1).Value = 45995
$ws.Cells.Item(48, 2).Value = 11
$ws.Cells.Item(48, 3).Value = 0.38
$ws.Cells.Item(48, 4).Value = "04.12.202511"
$ws.Cells.Item(49, 1).Value = 45995
$ws.Cells.Item(49, 2).Value = 12
$ws.Cells.Item(49, 3).Value = 0.931
$ws.Cells.Item(49, 4).Value = "04.12.202512"
$ws.Cells.Item(50, 1).Value = 45995
$ws.Cells.Item(50, 2).Value = 13
$ws.Cells.Item(50, 3).Value = 1.394
$ws.Cells.Item(50, 4).Value = "04.12.202513"
$ws.Cells.Item(51, 1).Value = 45995
$ws.Cells.Item(51, 2).Value = 14
$ws.Cells.Item(51, 3).Value = 1.577
$ws.Cells.Item(51, 4).Value = "04.12.202514"
$ws.Cells.Item(52, 1).Value = 45995
$ws.Cells.Item(52, 2).Value = 15
$ws.Cells.Item(52, 3).Value = 1.648
$ws.Cells.Item(52, 4).Value = "04.12.202515"
$ws.Cells.Item(53, 1).Value = 45995
$ws.Cells.Item(53, 2).Value = 16
$ws.Cells.Item(53, 3).Value = 0.98
$ws.Cells.Item(53, 4).Value = "04.12.202516"
$ws.Cells.Item(54, 1).Value = 45995
$ws.Cells.Item(54, 2).Value = 17
$ws.Cells.Item(54, 3).Value = 0.307
$ws.Cells.Item(54, 4).Value = "04.12.202517"
$ws.Cells.Item(55, 1).Value = 45995
$ws.Cells.Item(55, 2).Value = 18
$ws.Cells.Item(55, 3).Value = 0.022
$ws.Cells.Item(55, 4).Value = "04.12.202518"
$ws.Cells.Item(56, 1).Value = 45995
$ws.Cells.Item(56, 2).Value = 19
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(56, 4).Value = "04.12.202519"
$ws.Cells.Item(57, 1).Value = 45995
$ws.Cells.Item(57, 2).Value = 20
$ws.Cells.Item(57, 3).Value = 0
$ws.Cells.Item(57, 4).Value = "04.12.202520"
$ws.Cells.Item(58, 1).Value = 45995
$ws.Cells.Item(58, 2).Value = 21
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 4).Value = "04.12.202521"
$ws.Cells.Item(59, 1).Value = 45995
$ws.Cells.Item(59, 2).Value = 22
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = "04.12.202522"
$ws.Cells.Item(60, 1).Value = 45995
$ws.Cells.Item(60, 2).Value = 23
$ws.Cells.Item(60, 3).Value = 0
$ws.Cells.Item(60, 4).Value = "04.12.202523"
$ws.Cells.Item(61, 1).Value = 45995
$ws.Cells.Item(61, 2).Value = 24
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(61, 4).Value = "04.12.202524"
$ws.Cells.Item(62, 1).Value = 45996
$ws.Cells.Item(62, 2).Value = 1
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = "05.12.20251"
$ws.Cells.Item(63, 1).Value = 45996
$ws.Cells.Item(63, 2).Value = 2
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = "05.12.20252"
$ws.Cells.Item(64, 1).Value = 45996
$ws.Cells.Item(64, 2).Value = 3
$ws.Cells.Item(64, 3).Value = 0
$ws.Cells.Item(64, 4).Value = "05.12.20253"
$ws.Cells.Item(65, 1).Value = 45996
$ws.Cells.Item(65, 2).Value = 4
$ws.Cells.Item(65, 3).Value = 0
$ws.Cells.Item(65, 4).Value = "05.12.20254"
$ws.Cells.Item(66, 1).Value = 45996
$ws.Cells.Item(66, 2).Value = 5
$ws.Cells.Item(66, 3).Value = 0
$ws.Cells.Item(66, 4).Value = "05.12.20255"
$ws.Cells.Item(67, 1).Value = 45996
$ws.Cells.Item(67, 2).Value = 6
$ws.Cells.Item(67, 3).Value = 0
$ws.Cells.Item(67, 4).Value = "05.12.20256"
$ws.Cells.Item(68, 1).Value = 45996
$ws.Cells.Item(68, 2).Value = 7
$ws.Cells.Item(68, 3).Value = 0
$ws.Cells.Item(68, 4).Value = "05.12.20257"
$ws.Cells.Item(69, 1).Value = 45996
$ws.Cells.Item(69, 2).Value = 8
$ws.Cells.Item(69, 3).Value = 0
$ws.Cells.Item(69, 4).Value = "05.12.20258"
$ws.Cells.Item(70, 1).Value = 45996
$ws.Cells.Item(70, 2).Value = 9
$ws.Cells.Item(70, 3).Value = 0
$ws.Cells.Item(70, 4).Value = "05.12.20259"
$ws.Cells.Item(71, 1).Value = 45996
$ws.Cells.Item(71, 2).Value = 10
$ws.Cells.Item(71, 3).Value = 0.066
$ws.Cells.Item(71, 4).Value = "05.12.202510"
$ws.Cells.Item(72, 1).Value = 45996
$ws.Cells.Item(72, 2).Value = 11
$ws.Cells.Item(72, 3).Value = 0.348
$ws.Cells.Item(72, 4).Value = "05.12.202511"
$ws.Cells.Item(73, 1).Value = 45996
$ws.Cells.Item(73, 2).Value = 12
$ws.Cells.Item(73, 3).Value = 0.739
$ws.Cells.Item(73, 4).Value = "05.12.202512"
$ws.Cells.Item(74, 1).Value = 45996
$ws.Cells.Item(74, 2).Value = 13
$ws.Cells.Item(74, 3).Value = 1.006
$ws.Cells.Item(74, 4).Value = "05.12.202513"
$ws.Cells.Item(75, 1).Value = 45996
$ws.Cells.Item(75, 2).Value = 14
$ws.Cells.Item(75, 3).Value = 1.318
$ws.Cells.Item(75, 4).Value = "05.12.202514"
$ws.Cells.Item(76, 1).Value = 45996
$ws.Cells.Item(76, 2).Value = 15
$ws.Cells.Item(76, 3).Value = 1.093
$ws.Cells.Item(76, 4).Value = "05.12.202515"
$ws.Cells.Item(77, 1).Value = 45996
$ws.Cells.Item(77, 2).Value = 16
$ws.Cells.Item(77, 3).Value = 0.714
$ws.Cells.Item(77, 4).Value = "05.12.202516"
$ws.Cells.Item(78, 1).Value = 45996
$ws.Cells.Item(78, 2).Value = 17
$ws.Cells.Item(78, 3).Value = 0.225
$ws.Cells.Item(78, 4).Value = "05.12.202517"
$ws.Cells.Item(79, 1).Value = 45996
$ws.Cells.Item(79, 2).Value = 18
$ws.Cells.Item(79, 3).Value = 0.018
$ws.Cells.Item(79, 4).Value = "05.12.202518"
$ws.Cells.Item(80, 1).Value = 45996
$ws.Cells.Item(80, 2).Value = 19
$ws.Cells.Item(80, 3).Value = 0
$ws.Cells.Item(80, 4).Value = "05.12.202519"
$ws.Cells.Item(81, 1).Value = 45996
$ws.Cells.Item(81, 2).Value = 20
$ws.Cells.Item(81, 3).Value = 0
$ws.Cells.Item(81, 4).Value = "05.12.202520"
$ws.Cells.Item(82, 1).Value = 45996
$ws.Cells.Item(82, 2).Value = 21
$ws.Cells.Item(82, 3).Value = 0
$ws.Cells.Item(82, 4).Value = "05.12.202521"
$ws.Cells.Item(83, 1).Value = 45996
$ws.Cells.Item(83, 2).Value = 22
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = "05.12.202522"
$ws.Cells.Item(84, 1).Value = 45996
$ws.Cells.Item(84, 2).Value = 23
$ws.Cells.Item(84, 3).Value = 0
$ws.Cells.Item(84, 4).Value = "05.12.202523"
$ws.Cells.Item(85, 1).Value = 45996
$ws.Cells.Item(85, 2).Value = 24
$ws.Cells.Item(85, 3).Value = 0
$ws.Cells.Item(85, 4).Value = "05.12.202524"
$ws.Cells.Item(86, 1).Value = 45997
$ws.Cells.Item(86, 2).Value = 1
$ws.Cells.Item(86, 3).Value = 0
$ws.Cells.Item(86, 4).Value = "06.12.20251"
$ws.Cells.Item(87, 1).Value = 45997
$ws.Cells.Item(87, 2).Value = 2
$ws.Cells.Item(87, 3).Value = 0
$ws.Cells.Item(87, 4).Value = "06.12.20252"
$ws.Cells.Item(88, 1).Value = 45997
$ws.Cells.Item(88, 2).Value = 3
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = "06.12.20253"
$ws.Cells.Item(89, 1).Value = 45997
$ws.Cells.Item(89, 2).Value = 4
$ws.Cells.Item(89, 3).Value = 0
$ws.Cells.Item(89, 4).Value = "06.12.20254"
$ws.Cells.Item(90, 1).Value = 45997
$ws.Cells.Item(90, 2).Value = 5
$ws.Cells.Item(90, 3).Value = 0
$ws.Cells.Item(90, 4).Value = "06.12.20255"
$ws.Cells.Item(91, 1).Value = 45997
$ws.Cells.Item(91, 2).Value = 6
$ws.Cells.Item(91, 3).Value = 0
$ws.Cells.Item(91, 4).Value = "06.12.20256"
$ws.Cells.Item(92, 1).Value = 45997
$ws.Cells.Item(92, 2).Value = 7
$ws.Cells.Item(92, 3).Value = 0
$ws.Cells.Item(92, 4).Value = "06.12.20257"
$ws.Cells.Item(93, 1).Value = 45997
$ws.Cells.Item(93, 2).Value = 8
$ws.Cells.Item(93, 3).Value = 0
$ws.Cells.Item(93, 4).Value = "06.12.20258"
$ws.Cells.Item(94, 1).Value = 45997
$ws.Cells.Item(94, 2).Value = 9
$ws.Cells.Item(94, 3).Value = 0
$ws.Cells.Item(94, 4).Value = "06.12.20259"
$ws.Cells.Item(95, 1).Value = 45997
$ws.Cells.Item(95, 2).Value = 10
$ws.Cells.Item(95, 3).Value = 0.099
$ws.Cells.Item(95, 4).Value = "06.12.202510"
$ws.Cells.Item(96, 1).Value = 45997
$ws.Cells.Item(96, 2).Value = 11
$ws.Cells.Item(96, 3).Value = 0.588
$ws.Cells.Item(96, 4).Value = "06.12.202511"
$ws.Cells.Item(97, 1).Value = 45997
$ws.Cells.Item(97, 2).Value = 12
$ws.Cells.Item(97, 3).Value = 1.308
$ws.Cells.Item(97, 4).Value = "06.12.202512"
$ws.Cells.Item(98, 1).Value = 45997
$ws.Cells.Item(98, 2).Value = 13
$ws.Cells.Item(98, 3).Value = 1.736
$ws.Cells.Item(98, 4).Value = "06.12.202513"
$ws.Cells.Item(99, 1).Value = 45997
$ws.Cells.Item(99, 2).Value = 14
$ws.Cells.Item(99, 3).Value = 1.792
$ws.Cells.Item(99, 4).Value = "06.12.202514"
$ws.Cells.Item(100, 1).Value = 45997
$ws.Cells.Item(100, 2).Value = 15
$ws.Cells.Item(100, 3).Value = 1.655
$ws.Cells.Item(100, 4).Value = "06.12.202515"
$ws.Cells.Item(101, 1).Value = 45997
$ws.Cells.Item(101, 2).Value = 16
$ws.Cells.Item(101, 3).Value = 0.936
$ws.Cells.Item(101, 4).Value = "06.12.202516"
$ws.Cells.Item(102, 1).Value = 45997
$ws.Cells.Item(102, 2).Value = 17
$ws.Cells.Item(102, 3).Value = 0.238
$ws.Cells.Item(102, 4).Value = "06.12.202517"
$ws.Cells.Item(103, 1).Value = 45997
$ws.Cells.Item(103, 2).Value = 18
$ws.Cells.Item(103, 3).Value = 0.011
$ws.Cells.Item(103, 4).Value = "06.12.202518"
$ws.Cells.Item(104, 1).Value = 45997
$ws.Cells.Item(104, 2).Value = 19
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(104, 4).Value = "06.12.202519"
$ws.Cells.Item(105, 1).Value = 45997
$ws.Cells.Item(105, 2).Value = 20
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(105, 4).Value = "06.12.202520"
$ws.Cells.Item(106, 1).Value = 45997
$ws.Cells.Item(106, 2).Value = 21
$ws.Cells.Item(106, 3).Value = 0
$ws.Cells.Item(106, 4).Value = "06.12.202521"
$ws.Cells.Item(107, 1).Value = 45997
$ws.Cells.Item(107, 2).Value = 22
$ws.Cells.Item(107, 3).Value = 0
$ws.Cells.Item(107, 4).Value = "06.12.202522"
$ws.Cells.Item(108, 1).Value = 45997
$ws.Cells.Item(108, 2).Value = 23
$ws.Cells.Item(108, 3).Value = 0
$ws.Cells.Item(108, 4).Value = "06.12.202523"
$ws.Cells.Item(109, 1).Value = 45997
$ws.Cells.Item(109, 2).Value = 24
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = "06.12.202524"
$ws.Cells.Item(110, 1).Value = 45998
$ws.Cells.Item(110, 2).Value = 1
$ws.Cells.Item(110, 3).Value = 0
$ws.Cells.Item(110, 4).Value = "07.12.20251"
$ws.Cells.Item(111, 1).Value = 45998
$ws.Cells.Item(111, 2).Value = 2
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 4).Value = "07.12.20252"
$ws.Cells.Item(112, 1).Value = 45998
$ws.Cells.Item(112, 2).Value = 3
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = "07.12.20253"
$ws.Cells.Item(113, 1).Value = 45998
$ws.Cells.Item(113, 2).Value = 4
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 4).Value = "07.12.20254"
$ws.Cells.Item(114, 1).Value = 45998
$ws.Cells.Item(114, 2).Value = 5
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = "07.12.20255"
$ws.Cells.Item(115, 1).Value = 45998
$ws.Cells.Item(115, 2).Value = 6
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = "07.12.20256"
$ws.Cells.Item(116, 1).Value = 45998
$ws.Cells.Item(116, 2).Value = 7
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = "07.12.20257"
$ws.Cells.Item(117, 1).Value = 45998
$ws.Cells.Item(117, 2).Value = 8
$ws.Cells.Item(117, 3).Value = 0
$ws.Cells.Item(117, 4).Value = "07.12.20258"
$ws.Cells.Item(118, 1).Value = 45998
$ws.Cells.Item(118, 2).Value = 9
$ws.Cells.Item(118, 3).Value = 0
$ws.Cells.Item(118, 4).Value = "07.12.20259"
$ws.Cells.Item(119, 1).Value = 45998
$ws.Cells.Item(119, 2).Value = 10
$ws.Cells.Item(119, 3).Value = 0.067
$ws.Cells.Item(119, 4).Value = "07.12.202510"
$ws.Cells.Item(120, 1).Value = 45998
$ws.Cells.Item(120, 2).Value = 11
$ws.Cells.Item(120, 3).Value = 0.366
$ws.Cells.Item(120, 4).Value = "07.12.202511"
$ws.Cells.Item(121, 1).Value = 45998
$ws.Cells.Item(121, 2).Value = 12
$ws.Cells.Item(121, 3).Value = 0.816
$ws.Cells.Item(121, 4).Value = "07.12.202512"
$ws.Cells.Item(122, 1).Value = 45998
$ws.Cells.Item(122, 2).Value = 13
$ws.Cells.Item(122, 3).Value = 1.165
$ws.Cells.Item(122, 4).Value = "07.12.202513"
$ws.Cells.Item(123, 1).Value = 45998
$ws.Cells.Item(123, 2).Value = 14
$ws.Cells.Item(123, 3).Value = 1.225
$ws.Cells.Item(123, 4).Value = "07.12.202514"
$ws.Cells.Item(124, 1).Value = 45998
$ws.Cells.Item(124, 2).Value = 15
$ws.Cells.Item(124, 3).Value = 0.94
$ws.Cells.Item(124, 4).Value = "07.12.202515"
$ws.Cells.Item(125, 1).Value = 45998
$ws.Cells.Item(125, 2).Value = 16
$ws.Cells.Item(125, 3).Value = 0.636
$ws.Cells.Item(125, 4).Value = "07.12.202516"
$ws.Cells.Item(126, 1).Value = 45998
$ws.Cells.Item(126, 2).Value = 17
$ws.Cells.Item(126, 3).Value = 0.183
$ws.Cells.Item(126, 4).Value = "07.12.202517"
$ws.Cells.Item(127, 1).Value = 45998
$ws.Cells.Item(127, 2).Value = 18
$ws.Cells.Item(127, 3).Value = 0.01
$ws.Cells.Item(127, 4).Value = "07.12.202518"
$ws.Cells.Item(128, 1).Value = 45998
$ws.Cells.Item(128, 2).Value = 19
$ws.Cells.Item(128, 3).Value = 0
$ws.Cells.Item(128, 4).Value = "07.12.202519"
$ws.Cells.Item(129, 1).Value = 45998
$ws.Cells.Item(129, 2).Value = 20
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(129, 4).Value = "07.12.202520"
$ws.Cells.Item(130, 1).Value = 45998
$ws.Cells.Item(130, 2).Value = 21
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = "07.12.202521"
$ws.Cells.Item(131, 1).Value = 45998
$ws.Cells.Item(131, 2).Value = 22
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = "07.12.202522"
$ws.Cells.Item(132, 1).Value = 45998
$ws.Cells.Item(132, 2).Value = 23
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 4).Value = "07.12.202523"
$ws.Cells.Item(133, 1).Value = 45998
$ws.Cells.Item(133, 2).Value = 24
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(133, 4).Value = "07.12.202524"
$ws.Cells.Item(134, 1).Value = 45999
$ws.Cells.Item(134, 2).Value = 1
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(134, 4).Value = "08.12.20251"
$ws.Cells.Item(135, 1).Value = 45999
$ws.Cells.Item(135, 2).Value = 2
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 4).Value = "08.12.20252"
$ws.Cells.Item(136, 1).Value = 45999
$ws.Cells.Item(136, 2).Value = 3
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(136, 4).Value = "08.12.20253"
$ws.Cells.Item(137, 1).Value = 45999
$ws.Cells.Item(137, 2).Value = 4
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(137, 4).Value = "08.12.20254"
$ws.Cells.Item(138, 1).Value = 45999
$ws.Cells.Item(138, 2).Value = 5
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = "08.12.20255"
$ws.Cells.Item(139, 1).Value = 45999
$ws.Cells.Item(139, 2).Value = 6
$ws.Cells.Item(139, 3).Value = 0
$ws.Cells.Item(139, 4).Value = "08.12.20256"
$ws.Cells.Item(140, 1).Value = 45999
$ws.Cells.Item(140, 2).Value = 7
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(140, 4).Value = "08.12.20257"
$ws.Cells.Item(141, 1).Value = 45999
$ws.Cells.Item(141, 2).Value = 8
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(141, 4).Value = "08.12.20258"
$ws.Cells.Item(142, 1).Value = 45999
$ws.Cells.Item(142, 2).Value = 9
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(142, 4).Value = "08.12.20259"
$ws.Cells.Item(143, 1).Value = 45999
$ws.Cells.Item(143, 2).Value = 10
$ws.Cells.Item(143, 3).Value = 0.045
$ws.Cells.Item(143, 4).Value = "08.12.202510"
$ws.Cells.Item(144, 1).Value = 45999
$ws.Cells.Item(144, 2).Value = 11
$ws.Cells.Item(144, 3).Value = 0.286
$ws.Cells.Item(144, 4).Value = "08.12.202511"
$ws.Cells.Item(145, 1).Value = 45999
$ws.Cells.Item(145, 2).Value = 12
$ws.Cells.Item(145, 3).Value = 0.62
$ws.Cells.Item(145, 4).Value = "08.12.202512"
$ws.Cells.Item(146, 1).Value = 45999
$ws.Cells.Item(146, 2).Value = 13
$ws.Cells.Item(146, 3).Value = 0.79
$ws.Cells.Item(146, 4).Value = "08.12.202513"
$ws.Cells.Item(147, 1).Value = 45999
$ws.Cells.Item(147, 2).Value = 14
$ws.Cells.Item(147, 3).Value = 0.826
$ws.Cells.Item(147, 4).Value = "08.12.202514"
$ws.Cells.Item(148, 1).Value = 45999
$ws.Cells.Item(148, 2).Value = 15
$ws.Cells.Item(148, 3).Value = 0.717
$ws.Cells.Item(148, 4).Value = "08.12.202515"
$ws.Cells.Item(149, 1).Value = 45999
$ws.Cells.Item(149, 2).Value = 16
$ws.Cells.Item(149, 3).Value = 0.372
$ws.Cells.Item(149, 4).Value = "08.12.202516"
$ws.Cells.Item(150, 1).Value = 45999
$ws.Cells.Item(150, 2).Value = 17
$ws.Cells.Item(150, 3).Value = 0.116
$ws.Cells.Item(150, 4).Value = "08.12.202517"
$ws.Cells.Item(151, 1).Value = 45999
$ws.Cells.Item(151, 2).Value = 18
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = "08.12.202518"
$ws.Cells.Item(152, 1).Value = 45999
$ws.Cells.Item(152, 2).Value = 19
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 4).Value = "08.12.202519"
$ws.Cells.Item(153, 1).Value = 45999
$ws.Cells.Item(153, 2).Value = 20
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 4).Value = "08.12.202520"
$ws.Cells.Item(154, 1).Value = 45999
$ws.Cells.Item(154, 2).Value = 21
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(154, 4).Value = "08.12.202521"
$ws.Cells.Item(155, 1).Value = 45999
$ws.Cells.Item(155, 2).Value = 22
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = "08.12.202522"
$ws.Cells.Item(156, 1).Value = 45999
$ws.Cells.Item(156, 2).Value = 23
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = "08.12.202523"
$ws.Cells.Item(157, 1).Value = 45999
$ws.Cells.Item(157, 2).Value = 24
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = "08.12.202524"
$ws.Cells.Item(158, 1).Value = 46000
$ws.Cells.Item(158, 2).Value = 1
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = "09.12.20251"
$ws.Cells.Item(159, 1).Value = 46000
$ws.Cells.Item(159, 2).Value = 2
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(159, 4).Value = "09.12.20252"
$ws.Cells.Item(160, 1).Value = 46000
$ws.Cells.Item(160, 2).Value = 3
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = "09.12.20253"
$ws.Cells.Item(161, 1).Value = 46000
$ws.Cells.Item(161, 2).Value = 4
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = "09.12.20254"
$ws.Cells.Item(162, 1).Value = 46000
$ws.Cells.Item(162, 2).Value = 5
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = "09.12.20255"
$ws.Cells.Item(163, 1).Value = 46000
$ws.Cells.Item(163, 2).Value = 6
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(163, 4).Value = "09.12.20256"
$ws.Cells.Item(164, 1).Value = 46000
$ws.Cells.Item(164, 2).Value = 7
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = "09.12.20257"
$ws.Cells.Item(165, 1).Value = 46000
$ws.Cells.Item(165, 2).Value = 8
$ws.Cells.Item(165, 3).Value = 0
$ws.Cells.Item(165, 4).Value = "09.12.20258"
$ws.Cells.Item(166, 1).Value = 46000
$ws.Cells.Item(166, 2).Value = 9
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = "09.12.20259"
$ws.Cells.Item(167, 1).Value = 46000
$ws.Cells.Item(167, 2).Value = 10
$ws.Cells.Item(167, 3).Value = 0.025
$ws.Cells.Item(167, 4).Value = "09.12.202510"
$ws.Cells.Item(168, 1).Value = 46000
$ws.Cells.Item(168, 2).Value = 11
$ws.Cells.Item(168, 3).Value = 0.139
$ws.Cells.Item(168, 4).Value = "09.12.202511"
$ws.Cells.Item(169, 1).Value = 46000
$ws.Cells.Item(169, 2).Value = 12
$ws.Cells.Item(169, 3).Value = 0.336
$ws.Cells.Item(169, 4).Value = "09.12.202512"
$ws.Cells.Item(170, 1).Value = 46000
$ws.Cells.Item(170, 2).Value = 13
$ws.Cells.Item(170, 3).Value = 0.425
$ws.Cells.Item(170, 4).Value = "09.12.202513"
